$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.367.54"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "'1.567.08"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'0.9990"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'291.28"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "'0.3790"
$ws.Range("E7").Value = "  +3.16%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "'0.3406"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "'0.07610"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'1.138"
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "'5.984"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'6.936"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.563.94"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001133"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'89.93"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'0.06735"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'16.61"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").Value = "'22.348.56"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").Value = "'2.701"
$ws.Range("E26").Value = "  -7.37%  "
$ws.Range("D27").Value = "'20.09"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "'147.19"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").Value = "'5.019"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "'126.13"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").Value = "'1.745.63"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'2.015"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'6.092"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("D34").Value = "'0.9897"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").Value = "'10.10"
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.425"
$ws.Range("E36").Value = "  +7.40%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.08515"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").Value = "'0.02511"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "'0.2299"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "'0.06486"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'5.408"
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").Value = "'11.35"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").Value = "'0.6315"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "'0.9979"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'14.01"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").Value = "'3.810"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").Value = "'0.5926"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "'2.084"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").Value = "'1.252"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "'124.41"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "'0.07315"
$ws.Range("E51").Value = "  +0.45%  "
